$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title: "Testing" " " "custom" " " "properties" -> "Testing " "custom " "properties" ---
$title = $s.Shapes.Item(1)
$ttr = $title.TextFrame.TextRange

$ttr.Characters(1, 7).Text = "Testing "
$ttr.Characters(9, 1).Text = ""

$ttr.Characters(9, 6).Text = "custom "
$ttr.Characters(16, 1).Text = ""

# --- Subtitle: "This" " " "is" " " "a" " " "subtitle" <br> <br> "A." " " "M."
#     -> "This " "is " "a " "subtitle" <br> <br> "A. " "M." ---
$subtitle = $s.Shapes.Item(2)
$str = $subtitle.TextFrame.TextRange

$str.Characters(1, 4).Text = "This "
$str.Characters(6, 1).Text = ""

$str.Characters(6, 2).Text = "is "
$str.Characters(9, 1).Text = ""

$str.Characters(9, 1).Text = "a "
$str.Characters(11, 1).Text = ""

$str.Characters(21, 2).Text = "A. "
$str.Characters(24, 1).Text = ""
